$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.018541264468057
$ws.Range("D2").Value = 1.024128462187883
$ws.Range("E2").Value = 1.022155727825481
$ws.Range("F2").Value = 1.029979372652939
$ws.Range("I2").Value = 1.02905108162313
$ws.Range("J2").Value = 1.023749113728984
$ws.Range("K2").Value = 1.026957872351175
$ws.Range("L2").Value = 1.024990941089987
$ws.Range("M2").Value = 1.032791713615767
$ws.Range("N2").Value = 1.012003235603879
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.019467788996196
$ws.Range("D3").Value = 1.024795202116994
$ws.Range("E3").Value = 1.023028964857404
$ws.Range("F3").Value = 1.031140292148025
$ws.Range("I3").Value = 1.029229431520331
$ws.Range("J3").Value = 1.024312124823089
$ws.Range("K3").Value = 1.027432019784377
$ws.Range("L3").Value = 1.02567060191941
$ws.Range("M3").Value = 1.033759943817087
$ws.Range("N3").Value = 1.012190487044383
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.020067489029026
$ws.Range("D4").Value = 1.025226391496851
$ws.Range("E4").Value = 1.023594548949633
$ws.Range("F4").Value = 1.031891611754013
$ws.Range("I4").Value = 1.029343025660816
$ws.Range("J4").Value = 1.02467601809657
$ws.Range("K4").Value = 1.027737910382204
$ws.Range("L4").Value = 1.026110295955803
$ws.Range("M4").Value = 1.034386010737232
$ws.Range("N4").Value = 1.012311472664902
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.020319644572139
$ws.Range("D5").Value = 1.02540760581279
$ws.Range("E5").Value = 1.023832449071072
$ws.Range("F5").Value = 1.032207497040916
$ws.Range("I5").Value = 1.029390346722382
$ws.Range("J5").Value = 1.024828899366689
$ws.Range("K5").Value = 1.027866287015995
$ws.Range("L5").Value = 1.026295120482157
$ws.Range("M5").Value = 1.034649103188469
$ws.Range("N5").Value = 1.012362291993621
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.020361985020658
$ws.Range("D6").Value = 1.025438029074827
$ws.Range("E6").Value = 1.023872401057648
$ws.Range("F6").Value = 1.032260537363568
$ws.Range("I6").Value = 1.029398266675578
$ws.Range("J6").Value = 1.024854562951846
$ws.Range("K6").Value = 1.027887829105508
$ws.Range("L6").Value = 1.026326151969497
$ws.Range("M6").Value = 1.034693271345534
$ws.Range("N6").Value = 1.012370822251491
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.020070858180901
$ws.Range("D7").Value = 1.025228813119669
$ws.Range("E7").Value = 1.023597727278807
$ws.Range("F7").Value = 1.031895832508614
$ws.Range("I7").Value = 1.029343659672852
$ws.Range("J7").Value = 1.024678061295281
$ws.Range("K7").Value = 1.027739626621694
$ws.Range("L7").Value = 1.026112765680989
$ws.Range("M7").Value = 1.03438952660839
$ws.Range("N7").Value = 1.012312151884668
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.018854350788431
$ws.Range("D8").Value = 1.024353838192693
$ws.Range("E8").Value = 1.022450729663238
$ws.Range("F8").Value = 1.030371684586409
$ws.Range("I8").Value = 1.029111730279861
$ws.Range("J8").Value = 1.023939470968386
$ws.Range("K8").Value = 1.027118301544838
$ws.Range("L8").Value = 1.025220654194646
$ws.Range("M8").Value = 1.033119023173727
$ws.Range("N8").Value = 1.012066554917869
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.016712091237269
$ws.Range("D9").Value = 1.022810266786694
$ws.Range("E9").Value = 1.020433764413745
$ws.Range("F9").Value = 1.027686908928988
$ws.Range("I9").Value = 1.028689201591828
$ws.Range("J9").Value = 1.02263485118405
$ws.Range("K9").Value = 1.026016480152347
$ws.Range("L9").Value = 1.023647971665985
$ws.Range("M9").Value = 1.030876845707987
$ws.Range("N9").Value = 1.011632425893743
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.015284883607027
$ws.Range("D10").Value = 1.021780110426431
$ws.Range("E10").Value = 1.019091997923447
$ws.Range("F10").Value = 1.025897697840026
$ws.Range("I10").Value = 1.028398242655511
$ws.Range("J10").Value = 1.021763040291997
$ws.Range("K10").Value = 1.025277296391351
$ws.Range("L10").Value = 1.022599114835558
$ws.Range("M10").Value = 1.02937978335397
$ws.Range("N10").Value = 1.011342110634816
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.014667122878556
$ws.Range("D11").Value = 1.021333792792107
$ws.Range("E11").Value = 1.018511692572734
$ws.Range("F11").Value = 1.025123097670626
$ws.Range("I11").Value = 1.028270061120033
$ws.Range("J11").Value = 1.021385054801863
$ws.Range("K11").Value = 1.024956131766235
$ws.Range("L11").Value = 1.022144862459876
$ws.Range("M11").Value = 1.02873099768099
$ws.Range("N11").Value = 1.011216191579352
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.014437694009109
$ws.Range("D12").Value = 1.021167973635084
$ws.Range("E12").Value = 1.018296245661087
$ws.Range("F12").Value = 1.02483539727294
$ws.Range("I12").Value = 1.028222119564799
$ws.Range("J12").Value = 1.021244581781462
$ws.Range("K12").Value = 1.024836673415511
$ws.Range("L12").Value = 1.021976120023909
$ws.Range("M12").Value = 1.028489927281109
$ws.Range("N12").Value = 1.01116938822621
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.014486905690262
$ws.Range("D13").Value = 1.021203544059826
$ws.Range("E13").Value = 1.018342455019092
$ws.Range("F13").Value = 1.024897109060097
$ws.Range("I13").Value = 1.028232418081987
$ws.Range("J13").Value = 1.021274716987244
$ws.Range("K13").Value = 1.024862305024741
$ws.Range("L13").Value = 1.022012316420622
$ws.Range("M13").Value = 1.028541641437597
$ws.Range("N13").Value = 1.011179429121907
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.014648157500483
$ws.Range("D14").Value = 1.021320086875501
$ws.Range("E14").Value = 1.018493881530181
$ws.Range("F14").Value = 1.025099315848468
$ws.Range("I14").Value = 1.028266104969593
$ws.Range("J14").Value = 1.021373444731357
$ws.Range("K14").Value = 1.024946260633041
$ws.Range("L14").Value = 1.022130914416091
$ws.Range("M14").Value = 1.028711072411035
$ws.Range("N14").Value = 1.011212323438964
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.014747514710988
$ws.Range("D15").Value = 1.021391887892217
$ws.Range("E15").Value = 1.018587194253763
$ws.Range("F15").Value = 1.025223904840054
$ws.Range("I15").Value = 1.028286816960879
$ws.Range("J15").Value = 1.021434264579216
$ws.Range("K15").Value = 1.024997966809361
$ws.Range("L15").Value = 1.022203984870055
$ws.Range("M15").Value = 1.028815453485389
$ws.Range("N15").Value = 1.011232586558332
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.015325887230038
$ws.Range("D16").Value = 1.021809725820042
$ws.Range("E16").Value = 1.019130525459841
$ws.Range("F16").Value = 1.025949108414709
$ws.Range("I16").Value = 1.028406703447603
$ws.Range("J16").Value = 1.021788115762179
$ws.Range("K16").Value = 1.025298588060087
$ws.Range("L16").Value = 1.022629260224811
$ws.Range("M16").Value = 1.029422829562766
$ws.Range("N16").Value = 1.011350463050192
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.015688746836718
$ws.Range("D17").Value = 1.022071757531836
$ws.Range("E17").Value = 1.019471527717316
$ws.Range("F17").Value = 1.026404046636755
$ws.Range("I17").Value = 1.028481317895494
$ws.Range("J17").Value = 1.022009947656381
$ws.Range("K17").Value = 1.025486867669533
$ws.Range("L17").Value = 1.022896000506895
$ws.Range("M17").Value = 1.029803673719451
$ws.Range("N17").Value = 1.011424347633853
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.015900418826724
$ws.Range("D18").Value = 1.02222457169816
$ws.Range("E18").Value = 1.019670494924015
$ws.Range("F18").Value = 1.026669417946445
$ws.Range("I18").Value = 1.028524627482713
$ws.Range("J18").Value = 1.022139291553145
$ws.Range("K18").Value = 1.025596582455832
$ws.Range("L18").Value = 1.023051576826314
$ws.Range("M18").Value = 1.030025760832927
$ws.Range("N18").Value = 1.011467422922754
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.015972597242714
$ws.Range("D19").Value = 1.022276673190494
$ws.Range("E19").Value = 1.019738348827644
$ws.Range("F19").Value = 1.026759905015655
$ws.Range("I19").Value = 1.028539359000952
$ws.Range("J19").Value = 1.022183386507673
$ws.Range("K19").Value = 1.025633974441648
$ws.Range("L19").Value = 1.023104622820849
$ws.Range("M19").Value = 1.030101477843306
$ws.Range("N19").Value = 1.011482107026111
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.015649813143429
$ws.Range("D20").Value = 1.022043646534829
$ws.Range("E20").Value = 1.019434934531817
$ws.Range("F20").Value = 1.026355234673372
$ws.Range("I20").Value = 1.028473334371139
$ws.Range("J20").Value = 1.021986152039754
$ws.Range("K20").Value = 1.025466677958868
$ws.Range("L20").Value = 1.022867382701399
$ws.Range("M20").Value = 1.029762818212561
$ws.Range("N20").Value = 1.01141642261927
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.014600671892272
$ws.Range("D21").Value = 1.021285768931164
$ws.Range("E21").Value = 1.018449287309459
$ws.Range("F21").Value = 1.025039770424243
$ws.Range("I21").Value = 1.0282561941054
$ws.Range("J21").Value = 1.021344373864832
$ws.Range("K21").Value = 1.02492154230676
$ws.Range("L21").Value = 1.02209599061221
$ws.Range("M21").Value = 1.028661181491826
$ws.Range("N21").Value = 1.011202637747626
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.013941237714038
$ws.Range("D22").Value = 1.020809048167479
$ws.Range("E22").Value = 1.017830175871132
$ws.Range("F22").Value = 1.024212804770892
$ws.Range("I22").Value = 1.02811776482765
$ws.Range("J22").Value = 1.020940443972872
$ws.Range("K22").Value = 1.024577848229592
$ws.Range("L22").Value = 1.021610911630118
$ws.Range("M22").Value = 1.027968061254976
$ws.Range("N22").Value = 1.01106804110442
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.014290796888741
$ws.Range("D23").Value = 1.021061786687014
$ws.Range("E23").Value = 1.018158320994797
$ws.Range("F23").Value = 1.02465118378638
$ws.Range("I23").Value = 1.02819132918683
$ws.Range("J23").Value = 1.021154614325786
$ws.Range("K23").Value = 1.024760136358202
$ws.Range("L23").Value = 1.021868067969628
$ws.Range("M23").Value = 1.028335542741051
$ws.Range("N23").Value = 1.011139410467728
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.015667405535154
$ws.Range("D24").Value = 1.022056348760127
$ws.Range("E24").Value = 1.019451469212211
$ws.Range("F24").Value = 1.026377290656434
$ws.Range("I24").Value = 1.028476942436451
$ws.Range("J24").Value = 1.021996904399536
$ws.Range("K24").Value = 1.025475801146142
$ws.Range("L24").Value = 1.022880313882993
$ws.Range("M24").Value = 1.029781279220576
$ws.Range("N24").Value = 1.011420003655083
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.017265749173061
$ws.Range("D25").Value = 1.023209517544812
$ws.Range("E25").Value = 1.020954695499509
$ws.Range("F25").Value = 1.028380874261599
$ws.Range("I25").Value = 1.028800072147627
$ws.Range("J25").Value = 1.022972493411607
$ws.Range("K25").Value = 1.026302148055789
$ws.Range("L25").Value = 1.024054622072892
$ws.Range("M25").Value = 1.03145690375903
$ws.Range("N25").Value = 1.011744817702694
